$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix company name text in B3 (PLC -> Plc)
$ws.Range("B3").Value = "Dar es Salaam Stock Exchange Plc (DAR:DSE)"

# Update numeric columns for row 2 and row 3 (same values for both rows)
$cols = @{
    "G"  = 0.3346303501945526
    "H"  = 0.3346303501945526
    "I"  = 0.166147859922179
    "J"  = 0.1654094249891915
    "K"  = 1.34
    "L"  = 0.5214007782101168
    "M"  = 0.765
    "N"  = 0.08462389380530974
    "O"  = 0.5708955223880596
    "P"  = 0.765
    "Q"  = 0.08462389380530974
    "R"  = 0.5708955223880596
    "U"  = 0.171
    "V"  = 0.01891592920353983
    "W"  = 0.1420996818663839
    "X"  = 0.02580261134732567
    "Y"  = 0.1162970705190582
    "Z"  = 0.2795605351898183
    "AA" = 0.04624194737541849
    "AB" = 0.02580261134732567
    "AC" = 0.02043933602809282
    "AG" = -0.171
    "AJ" = -0.01928064043296877
    "AK" = -0.01768538628606888
    "AM" = -0.707
    "AP" = -0.3263358778625954
    "AQ" = -0.6039603960396039
}

foreach ($col in $cols.Keys) {
    $value = $cols[$col]
    $ws.Range("$col`2").Value = $value
    $ws.Range("$col`3").Value = $value
}
